# Applies the commit "Se corrigieron faltas de ortografía y se hicieron
# unos ultimos arreglos":
#   1. Refresh the cached "datetimeFigureOut" footer field from 01/10/2019
#      to 02/10/2019 on the slide master and every slide layout.
#   2. Slide 8: drop the trailing "($/ton)" from the profit sentence.
#   3. Slide 9: merge the "Las " run into the following run so the whole
#      sentence lives in a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder refresh (master + every custom layout)
# ---------------------------------------------------------------------
function Update-DatePlaceholders($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $len = $tr.Length
            if ($len -gt 0) {
                $whole = $tr.Characters(1, $len)
                $whole.Text = $newText
            } else {
                $tr.Text = $newText
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes "02/10/2019"
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $cl = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $cl.Shapes "02/10/2019"
}

# ---------------------------------------------------------------------
# 2) Slide 8 - remove "($/ton)" from the profit sentence
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(2)
$tr8 = $shp8.TextFrame.TextRange
$para8 = $tr8.Paragraphs(3)
$len8 = $para8.Length
$whole8 = $para8.Characters(1, $len8)
$whole8.Text = "Fabricando alimentos con los siguientes componentes.243 de VEG2 y .756 de OIL1, tendríamos una ganancia de 127.56."

# ---------------------------------------------------------------------
# 3) Slide 9 - merge "Las " into the following run
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(2)
$tr9 = $shp9.TextFrame.TextRange
$para9 = $tr9.Paragraphs(3)
$prefix9 = $para9.Characters(1, 4)
$prefix9.Text = ""
$len9 = $para9.Length
$whole9 = $para9.Characters(1, $len9)
$whole9.Text = "Las ganancias máximas durante los primeros 6 meses del año no son tan buenas."
